# Update specific imputed values in the RandomForest result sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = 6.148299999999997
$ws.Range("A3").Value  = -21.49410000000002
$ws.Range("B5").Value  = 5.2986
$ws.Range("D5").Value  = -8.847599999999996
$ws.Range("D9").Value  = -8.571800000000003
$ws.Range("D11").Value = -8.229000000000003
$ws.Range("A14").Value = -20.70979999999999
$ws.Range("A16").Value = -21.18240000000002
$ws.Range("B16").Value = 6.208499999999995
$ws.Range("D17").Value = -8.539600000000005
$ws.Range("A21").Value = -21.35380000000001
$ws.Range("D21").Value = -8.182400000000007
$ws.Range("A23").Value = -21.36290000000003
$ws.Range("A25").Value = -22.33270000000003

$wb.Save()
